$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Row 2
$ws.Range("B2").Value = "<do>"
$ws.Range("C2").Value = 36

# Row 3
$ws.Range("B3").Value = "<kilo>"
$ws.Range("C3").Value = 34

# Row 4
$ws.Range("B4").Value = "<now>"
$ws.Range("C4").Value = 32

# Row 5
$ws.Range("B5").Value = "<echse>"
$ws.Range("C5").Value = 33

# Row 7
$ws.Range("B7").Value = "<line>"

# Row 8
$ws.Range("B8").Value = "<for>"
$ws.Range("C8").Value = 35

# Row 10
$ws.Range("B10").Value = "<all>"
$ws.Range("C10").Value = 35

# Row 11
$ws.Range("B11").Value = "<an>"
$ws.Range("C11").Value = 32

# Row 12
$ws.Range("C12").Value = 32

# Row 14
$ws.Range("B14").Value = "<six>"
$ws.Range("C14").Value = 35

# Row 15
$ws.Range("C15").Value = 31

# Row 16
$ws.Range("C16").Value = 35

# Row 17
$ws.Range("B17").Value = "<there>"
$ws.Range("C17").Value = 28

# Row 18
$ws.Range("C18").Value = 29
